$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: a new measurement that starts and ends on a different day (J7)
$ws.Range("A8").Value = 44348
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "23:58:30"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "00:02:18"
$ws.Range("D8").Value = "J7"
$ws.Range("E8").Value = "0.4"
$ws.Range("F8").Value = 17.2
$ws.Range("H8").Value = "Floating short"

# Update current selection to C9
$ws.Range("C9").Select()
